$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the dataset. It belongs right before the
# existing row 323 (chronologically among the other records in that block),
# so insert a new row there which pushes all the following rows down by one.
$ws.Rows.Item(323).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(323, 1).Value = 5
$ws.Cells.Item(323, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(323, 3).Value = "Maule"
$ws.Cells.Item(323, 4).Value = 45034
$ws.Cells.Item(323, 5).Value = 7
$ws.Cells.Item(323, 6).Value = "Fruta"
$ws.Cells.Item(323, 7).Value = 100108
$ws.Cells.Item(323, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(323, 9).Value = 100108005
$ws.Cells.Item(323, 10).Value = "Piña"
$ws.Cells.Item(323, 11).Value = "Caramelo"
$ws.Cells.Item(323, 12).Value = "Tercera"
$ws.Cells.Item(323, 13).Value = 250
$ws.Cells.Item(323, 14).Value = 20000
$ws.Cells.Item(323, 15).Value = 20000
$ws.Cells.Item(323, 16).Value = 20000
$ws.Cells.Item(323, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(323, 18).Value = "Ecuador"
$ws.Cells.Item(323, 19).Value = 1250
$ws.Cells.Item(323, 20).Value = 16
